$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values in columns A and B (rows 1-32)
$ws.Range("A1").Value = -0.34985991912974157
$ws.Range("B1").Value = 0.34892506166082171
$ws.Range("A2").Value = -0.2207104345619868
$ws.Range("B2").Value = 0.21869507509484265
$ws.Range("A3").Value = -0.11574613204672346
$ws.Range("B3").Value = 0.11529328227173608
$ws.Range("A4").Value = -0.10329328233475366
$ws.Range("B4").Value = 0.10291471868270108
$ws.Range("A5").Value = -0.096914718918728049
$ws.Range("B5").Value = 0.096184069002950068
$ws.Range("A6").Value = -0.037667284197057338
$ws.Range("B6").Value = 0.037646140446406573
$ws.Range("A7").Value = -0.01764614073357329
$ws.Range("B7").Value = 0.017629596614574083
$ws.Range("A8").Value = -0.018022595488451643
$ws.Range("B8").Value = 0.01799394603444604
$ws.Range("A9").Value = -0.01199394627723116
$ws.Range("B9").Value = 0.011979821731149976
$ws.Range("A10").Value = -0.0059798219745346159
$ws.Range("B10").Value = 0.0059802937632227327
$ws.Range("A11").Value = -0.0014802940017268895
$ws.Range("B11").Value = 0.0014805612759438702
$ws.Range("A12").Value = 0.0045194384806905852
$ws.Range("B12").Value = -0.0045214769137205124
$ws.Range("A13").Value = 0.010521476670800389
$ws.Range("B13").Value = -0.010525894989088513
$ws.Range("A14").Value = 0.022525894726848072
$ws.Range("B14").Value = -0.022551198581009402
$ws.Range("A15").Value = -0.021051742617899549
$ws.Range("B15").Value = 0.021026961284206003
$ws.Range("A16").Value = -0.01502696152687033
$ws.Range("B16").Value = 0.015004454734373951
$ws.Range("A17").Value = -0.009004454978104981
$ws.Range("B17").Value = 0.0089999997461829295
$ws.Range("A18").Value = -0.036111471373907023
$ws.Range("B18").Value = 0.036097161072788708
$ws.Range("A19").Value = -0.027097161308671236
$ws.Range("B19").Value = 0.027013976119574235
$ws.Range("A20").Value = -0.018013976357673656
$ws.Range("B20").Value = 0.018004309348860659
$ws.Range("A21").Value = -0.0090043095873202361
$ws.Range("B21").Value = 0.0089999997612908444
$ws.Range("A22").Value = -0.093948943289600351
$ws.Range("B22").Value = 0.093635229891336991
$ws.Range("A23").Value = -0.084635230135147843
$ws.Range("B23").Value = 0.084127041076286879
$ws.Range("A24").Value = -0.0421270414324173
$ws.Range("B24").Value = 0.041999999641814334
$ws.Range("A25").Value = -0.05251678504679802
$ws.Range("B25").Value = 0.052453592782917724
$ws.Range("A26").Value = -0.046453593024313733
$ws.Range("B26").Value = 0.04637888134767465
$ws.Range("A27").Value = -0.040378881589721694
$ws.Range("B27").Value = 0.040145923855770871
$ws.Range("A28").Value = -0.034145924100413616
$ws.Range("B28").Value = 0.034001720316335415
$ws.Range("A29").Value = -0.022001720582313311
$ws.Range("B29").Value = 0.021951266725411145
$ws.Range("A30").Value = -0.0019512670182044722
$ws.Range("B30").Value = 0.0018989596221690164
$ws.Range("A31").Value = -0.027019632242065938
$ws.Range("B31").Value = 0.027000904052494334
$ws.Range("A32").Value = -0.006000904349096281
$ws.Range("B32").Value = 0.0059999997524515791

# Widen column B to match column A's width (15.42578125 characters)
$ws.Columns.Item(2).ColumnWidth = 15.42578125
